$wb = $excel.ActiveWorkbook

# --- Sheet "LIST" (sheet1) ---
$wsList = $wb.Worksheets.Item("LIST")

# Update A3: RO.ORG.001.CRE -> RO.ORG.001.LEC
$wsList.Range("A3").Value = "RO.ORG.001.LEC"

# Update A4: RO.ORG.001.LEC -> RO.ACT.001.LEC
$wsList.Range("A4").Value = "RO.ACT.001.LEC"

# Add new cell C4: RO.ORG.001.LEC
$wsList.Range("C4").Value = "RO.ORG.001.LEC"

# Move selection to A4
$wsList.Activate()
$wsList.Range("A4").Select()

# --- Sheet "Feuil1" (sheet2) ---
$wsFeuil = $wb.Worksheets.Item("Feuil1")

# Move selection to C4
$wsFeuil.Activate()
$wsFeuil.Range("C4").Select()

# Reactivate LIST sheet as the originally tab-selected sheet
$wsList.Activate()
